# "Otočené poradie v tabuľkách (najnovšie vždy hore)."
# Adds a new competition results sheet "asdas" (1. Výročie FRI) after the
# existing sheets, so the newest table sits last / on top of the tab order.

$wb = $excel.ActiveWorkbook

# Reference sheet used only to source an already-blank A1 cell (keeps the
# new sheet's corner cell consistent with the other result sheets).
$refSheet = $wb.Worksheets.Item(1)

# Add the new worksheet after the last existing one, so the tab order stays
# "Testovacia súťaž", "Ukončená súťaž", "asdas".
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "asdas"

# A1 mirrors the blank corner cell used on the other result sheets.
$refSheet.Range("A1").Copy($newSheet.Range("A1"))

# Header row.
$newSheet.Cells.Item(1, 2).Value = "1. Výročie FRI"
$newSheet.Cells.Item(1, 4).Value = "Počet správnych odpovedí"
$newSheet.Cells.Item(1, 5).Value = "Úspešnosť (v %)"

# Single response row.
$newSheet.Cells.Item(3, 1).Value = "admin@frivia.sk"
$newSheet.Cells.Item(3, 2).Value = "Áno"
$newSheet.Cells.Item(3, 4).Value = 1
$newSheet.Cells.Item(3, 5).Value = 100
